$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row below row 13 (shifts rows 14:102 down to 15:103)
$ws.Rows.Item(14).Insert()

# Update B13's total_venda for day 12 of May/2025
$ws.Range("B13").Value = 31420.44

# New row 14 becomes day 13 of May/2025
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 26203.86
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 2025
$ws.Range("E14").Value = "05/2025"
